$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "operations provided by the " -> "operations provided by "
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "operations provided by the ", $true, $false, $false, $false, $false,
    $true, 1, $false, "operations provided by ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. "Common data engineering tasks " -> "Performing data engineering
#    tasks like " and insert "and " before "need to sum a large set..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Common data engineering tasks ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Performing data engineering tasks like ", 2) | Out-Null

$d.Content.Find.Execute(
    "need to sum a large set of numbers", $true, $false, $false, $false, $false,
    $true, 1, $false, "and need to sum a large set of numbers", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. "This paper is of interest" -> "This article is of interest"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "This paper is of interest", $true, $false, $false, $false, $false,
    $true, 1, $false, "This article is of interest", 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Remove the blank paragraph just before the "Theoretical Background"
#    Heading2 paragraph.
# ---------------------------------------------------------------------
$findHeading = $d.Content
$findHeading.Find.Execute(
    "Theoretical Background", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$headingPara = $findHeading.Paragraphs(1)
$headingPara.Previous().Range.Delete()

# ---------------------------------------------------------------------
# 5. "The resources section below offers some numerical analysis
#    articles to help understand the theory." ->
#    "The *resources section* below offers some numerical analysis
#    articles to help understand the theory behind the problem of
#    errors in summation."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The resources section below offers some numerical analysis articles to help understand the theory.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The resources section below offers some numerical analysis articles to help understand the theory behind the problem of errors in summation.",
    2) | Out-Null

$italicTarget = $d.Content
$italicTarget.Find.Execute(
    "resources section", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$italicTarget.Font.Italic = $true

# ---------------------------------------------------------------------
# 6. Heading2 "The Addition Problem in Standard Programming" gets
#    w:before="0" added to its paragraph spacing.
# ---------------------------------------------------------------------
$headingSpacing = $d.Content
$headingSpacing.Find.Execute(
    "The Addition Problem in", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$headingSpacing.ParagraphFormat.SpaceBefore = 0

# ---------------------------------------------------------------------
# 7. "We introduce representation error by dividing sequences by a
#    large prime number" ->
#    "We next introduce representation error by dividing the integer
#    sequences by a large prime number"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "We introduce representation error ", $true, $false, $false, $false, $false,
    $true, 1, $false, "We next introduce representation error ", 2) | Out-Null

$d.Content.Find.Execute(
    "by dividing sequences by a large prime number", $true, $false, $false, $false, $false,
    $true, 1, $false, "by dividing the integer sequences by a large prime number", 2) | Out-Null

# ---------------------------------------------------------------------
# 8. "summation algorithm " + "(see reference 5 below.)" -> merge text
#    "summation algorithm (see reference 5 below.)"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "summation algorithm (see reference 5 below.)", $true, $false, $false, $false, $false,
    $true, 1, $false, "summation algorithm (see reference 5 below.)", 2) | Out-Null

# ---------------------------------------------------------------------
# 9. "Experiment Summary:" -> "Experiment Summary"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Experiment Summary:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Experiment Summary", 2) | Out-Null

# ---------------------------------------------------------------------
# 10. ", yielding 7 subsequences" -> ", yielding 7 sub-sequences"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    ", yielding 7 subsequences", $true, $false, $false, $false, $false,
    $true, 1, $false, ", yielding 7 sub-sequences", 2) | Out-Null

# ---------------------------------------------------------------------
# 11. "This graph shows we lose nearly four" -> "This shows we lose
#     nearly four" and move the _GoBack bookmark here, between "This "
#     and "shows".
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "This graph shows we lose nearly four", $true, $false, $false, $false, $false,
    $true, 1, $false, "This shows we lose nearly four", 2) | Out-Null

$bookmarkAnchor = $d.Content
$bookmarkAnchor.Find.Execute(
    "This shows we lose nearly four", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$bookmarkPos = $bookmarkAnchor.Start + 5
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------
# 12. "A good Finite Precision Math Tutorial is" + " found at " -> merge
#     text "A good Finite Precision Math Tutorial is found at "
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "A good Finite Precision Math Tutorial is found at ", $true, $false, $false, $false, $false,
    $true, 1, $false, "A good Finite Precision Math Tutorial is found at ", 2) | Out-Null
